# Powerpoint writer: consolidate text run nodes.
# Merge each "word" run with its immediately following single-space run
# into one run (e.g. "Testing" + " " -> "Testing "), matching the
# consolidated <a:r> layout produced by the updated writer.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1 (Title): "Testing" " " "custom" " " "properties" ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Characters(1, 8).Text = "Testing "   # "Testing" + " " -> "Testing "
$title.Characters(9, 7).Text = "custom "    # "custom" + " " -> "custom "

# --- Shape 2 (Subtitle): "This" " " "is" " " "a" " " "subtitle" <br/><br/> "A." " " "M." ---
$subtitle = $s.Shapes.Item(2).TextFrame.TextRange
$subtitle.Characters(1, 5).Text = "This "   # "This" + " " -> "This "
$subtitle.Characters(6, 3).Text = "is "     # "is" + " " -> "is "
$subtitle.Characters(9, 2).Text = "a "      # "a" + " " -> "a "
$subtitle.Characters(21, 3).Text = "A. "    # "A." + " " -> "A. "
